# Add data for 2022-06-20 (through June 12 snapshot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab and update the matching header label in B1's text.
$ws.Name = "Through 2022-06-12"
$ws.Range("B1").Value = "June 2022 (through June 12)"

# Row 3 - Auburn Gresham
$ws.Range("B3").Value = 3
$ws.Range("H3").Value = 1

# Row 4 - North Lawndale
$ws.Range("H4").Value = 1
$ws.Range("N4").Value = 2
$ws.Range("Z4").Value = 4

# Row 10 - Garfield Park
$ws.Range("Z10").Value = 1

# Row 12 - Roseland
$ws.Range("H12").Value = 1

# Row 14 - Austin
$ws.Range("B14").Value = 2
$ws.Range("H14").Value = 4
$ws.Range("N14").Value = 3

# Row 18 - Lake View
$ws.Range("H18").Value = 1

# Row 19 - Little Italy, UIC
$ws.Range("T19").Value = 2

# Row 20 - Hyde Park
$ws.Range("AF20").Value = 1

# Row 25 - Ashburn
$ws.Range("Z25").Value = 1

# Row 55 - East Village
$ws.Range("B55").Value = 1

# Row 70 - Loop
$ws.Range("B70").Value = 2

# Row 71 - Lower West Side
$ws.Range("AR71").Value = 1
